$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 440762.72
$ws.Range("J17").Value = 440762.72
$ws.Range("L17").Value = 1322288.16
$ws.Range("N17").Value = -1322624.16

$ws.Range("H76").Value = 60002760
$ws.Range("I76").Value = 62502750
$ws.Range("K76").Value = 62502750
$ws.Range("M76").Value = -62502435

$ws.Range("H79").Value = 60002760
$ws.Range("I79").Value = 62502750
$ws.Range("K79").Value = 62502750
$ws.Range("M79").Value = -62501658

$ws.Range("H98").Value = 1486.7812
$ws.Range("I98").Value = 1496.0322
$ws.Range("K98").Value = 1496.0322
$ws.Range("M98").Value = 1.967799999999897

$ws.Range("H112").Value = 2941.3794
$ws.Range("I112").Value = 1640
$ws.Range("J112").Value = 3212.5
$ws.Range("K112").Value = 4920
$ws.Range("L112").Value = 9637.5
$ws.Range("M112").Value = -3812
$ws.Range("N112").Value = -11853.5

$ws.Range("H122").Value = 1486.7812
$ws.Range("I122").Value = 1496.0322
$ws.Range("K122").Value = 4488.096600000001
$ws.Range("M122").Value = -2038.096600000001

$ws.Range("H129").Value = 727210.9
$ws.Range("J129").Value = 904477
$ws.Range("L129").Value = 2713431
$ws.Range("N129").Value = -2723431

$ws.Range("H132").Value = 1756193.1
$ws.Range("I132").Value = 1710.9348
$ws.Range("K132").Value = 5132.8044
$ws.Range("M132").Value = -2602.8044

$ws.Range("H137").Value = 1160.4
$ws.Range("I137").Value = 960.4
$ws.Range("K137").Value = 2881.2
$ws.Range("M137").Value = -331.1999999999998

$ws.Range("H138").Value = 3045.6
$ws.Range("I138").Value = 1375.4584
$ws.Range("J138").Value = 3573.0132
$ws.Range("K138").Value = 4126.3752
$ws.Range("L138").Value = 10719.0396
$ws.Range("M138").Value = 1013.6248
$ws.Range("N138").Value = -20999.0396

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 225.5
$ws.Range("I4").Value = 133.33333
$ws.Range("J4").Value = 502
$ws.Range("K4").Value = 133.33333
$ws.Range("L4").Value = 502
$ws.Range("M4").Value = -17.33332999999999
$ws.Range("N4").Value = -734

$ws.Range("H32").Value = 17361.482
$ws.Range("I32").Value = 12759.68
$ws.Range("K32").Value = 12759.68
$ws.Range("M32").Value = -12472.68

$ws.Range("H63").Value = 2057
$ws.Range("I63").Value = 2004.2307
$ws.Range("K63").Value = 2004.2307
$ws.Range("M63").Value = -1318.2307

$ws.Range("H66").Value = 2057
$ws.Range("I66").Value = 2004.2307
$ws.Range("K66").Value = 10021.1535
$ws.Range("M66").Value = -6589.1535

$ws.Range("H97").Value = 1894.4872
$ws.Range("I97").Value = 2491.16
$ws.Range("J97").Value = 829
$ws.Range("K97").Value = 2491.16
$ws.Range("L97").Value = 829
$ws.Range("M97").Value = -1995.16
$ws.Range("N97").Value = -1821

$ws.Range("H132").Value = 1645.7637
$ws.Range("I132").Value = 1247.9
$ws.Range("K132").Value = 3743.7
$ws.Range("M132").Value = -1213.7

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19000
$ws.Range("J82").Value = 36000
$ws.Range("L82").Value = 36000
$ws.Range("N82").Value = -36766

$ws.Range("H85").Value = 19000
$ws.Range("J85").Value = 36000
$ws.Range("L85").Value = 36000
$ws.Range("N85").Value = -38652

$ws.Range("H99").Value = 34485020
$ws.Range("I99").Value = 38463600
$ws.Range("J99").Value = 3993.6667
$ws.Range("K99").Value = 38463600
$ws.Range("L99").Value = 3993.6667
$ws.Range("M99").Value = -38462102
$ws.Range("N99").Value = -6989.6667

$ws.Range("H134").Value = 26626.244
$ws.Range("I134").Value = 1918.9429
$ws.Range("J134").Value = 170752.17
$ws.Range("K134").Value = 5756.8287
$ws.Range("L134").Value = 512256.51
$ws.Range("M134").Value = -3221.8287
$ws.Range("N134").Value = -517326.51

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 114
$ws.Range("I7").Value = 83
$ws.Range("J7").Value = 300
$ws.Range("K7").Value = 83
$ws.Range("L7").Value = 300
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = -526

$ws.Range("H31").Value = 2334.16
$ws.Range("I31").Value = 1918.9025
$ws.Range("J31").Value = 4225.8887
$ws.Range("K31").Value = 1918.9025
$ws.Range("L31").Value = 4225.8887
$ws.Range("M31").Value = -1623.9025
$ws.Range("N31").Value = -4815.8887

$ws.Range("H34").Value = 2334.16
$ws.Range("I34").Value = 1918.9025
$ws.Range("J34").Value = 4225.8887
$ws.Range("K34").Value = 1918.9025
$ws.Range("L34").Value = 4225.8887
$ws.Range("M34").Value = -1716.9025
$ws.Range("N34").Value = -4629.8887

$ws.Range("H58").Value = 2700.6038
$ws.Range("I58").Value = 683.4186
$ws.Range("K58").Value = 683.4186
$ws.Range("M58").Value = -480.4186

$ws.Range("H136").Value = 2700.6038
$ws.Range("I136").Value = 683.4186
$ws.Range("K136").Value = 2050.2558
$ws.Range("M136").Value = 499.7442000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 557.625
$ws.Range("J34").Value = 600
$ws.Range("L34").Value = 1800
$ws.Range("N34").Value = -1968

$ws.Range("H39").Value = 6413856
$ws.Range("J39").Value = 6413856
$ws.Range("L39").Value = 19241568
$ws.Range("N39").Value = -19242156

$ws.Range("H55").Value = 125975.125
$ws.Range("I55").Value = 1000000
$ws.Range("J55").Value = 1114.4286
$ws.Range("K55").Value = 3000000
$ws.Range("L55").Value = 3343.2858
$ws.Range("M55").Value = -2999823
$ws.Range("N55").Value = -3697.2858

$ws.Range("H69").Value = 387.5
$ws.Range("J69").Value = 387.5
$ws.Range("L69").Value = 1162.5
$ws.Range("N69").Value = -2784.5

$ws.Range("H72").Value = 387.5
$ws.Range("J72").Value = 387.5
$ws.Range("L72").Value = 3487.5
$ws.Range("N72").Value = -11599.5

$ws.Range("H115").Value = 1030
$ws.Range("I115").Value = 450
$ws.Range("J115").Value = 1465
$ws.Range("K115").Value = 1350
$ws.Range("L115").Value = 4395
$ws.Range("M115").Value = -175
$ws.Range("N115").Value = -6745

$ws.Range("H122").Value = 21186.28
$ws.Range("I122").Value = 517.2857
$ws.Range("J122").Value = 24551
$ws.Range("K122").Value = 4655.571300000001
$ws.Range("L122").Value = 220959
$ws.Range("M122").Value = -2205.571300000001
$ws.Range("N122").Value = -225859

$ws.Range("H131").Value = 70268.555
$ws.Range("I131").Value = 83793.336
$ws.Range("J131").Value = 60721.65
$ws.Range("K131").Value = 251380.008
$ws.Range("L131").Value = 182164.95
$ws.Range("M131").Value = -246340.008
$ws.Range("N131").Value = -192244.95

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 7387.857
$ws.Range("J33").Value = 7387.857
$ws.Range("L33").Value = 7387.857
$ws.Range("N33").Value = -7891.857

$ws.Range("H40").Value = 10996
$ws.Range("J40").Value = 10996
$ws.Range("L40").Value = 10996
$ws.Range("N40").Value = -11298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 65755
$ws.Range("J46").Value = 65755
$ws.Range("L46").Value = 65755
$ws.Range("N46").Value = -66217

$ws.Range("H96").Value = 2586.2856
$ws.Range("I96").Value = 3150
$ws.Range("J96").Value = 2360.8
$ws.Range("K96").Value = 3150
$ws.Range("L96").Value = 2360.8
$ws.Range("M96").Value = -1777
$ws.Range("N96").Value = -5106.8

$ws.Range("H134").Value = 65755
$ws.Range("J134").Value = 65755
$ws.Range("L134").Value = 197265
$ws.Range("N134").Value = -202335
